$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Cells.Item(17, 8).Value = 1833.9535
$ws.Cells.Item(17, 10).Value = 2026.0605
$ws.Cells.Item(17, 12).Value = 6078.181500000001
$ws.Cells.Item(17, 14).Value = -6414.181500000001

# Row 53
$ws.Cells.Item(53, 8).Value = 407.0909
$ws.Cells.Item(53, 10).Value = 525.7143
$ws.Cells.Item(53, 12).Value = 525.7143
$ws.Cells.Item(53, 14).Value = -1799.7143

# Row 62
$ws.Cells.Item(62, 8).Value = 13432.258
$ws.Cells.Item(62, 9).Value = 13771.421
$ws.Cells.Item(62, 10).Value = 12895.25
$ws.Cells.Item(62, 11).Value = 13771.421
$ws.Cells.Item(62, 12).Value = 12895.25
$ws.Cells.Item(62, 13).Value = -13147.421
$ws.Cells.Item(62, 14).Value = -14143.25

# Row 65
$ws.Cells.Item(65, 8).Value = 13432.258
$ws.Cells.Item(65, 9).Value = 13771.421
$ws.Cells.Item(65, 10).Value = 12895.25
$ws.Cells.Item(65, 11).Value = 68857.105
$ws.Cells.Item(65, 12).Value = 64476.25
$ws.Cells.Item(65, 13).Value = -65737.105
$ws.Cells.Item(65, 14).Value = -70716.25

# Row 70
$ws.Cells.Item(70, 8).Value = 21328.5
$ws.Cells.Item(70, 9).Value = 1992
$ws.Cells.Item(70, 10).Value = 60001.5
$ws.Cells.Item(70, 11).Value = 5976
$ws.Cells.Item(70, 12).Value = 180004.5
$ws.Cells.Item(70, 13).Value = -5706
$ws.Cells.Item(70, 14).Value = -180544.5

# Row 73
$ws.Cells.Item(73, 8).Value = 21328.5
$ws.Cells.Item(73, 9).Value = 1992
$ws.Cells.Item(73, 10).Value = 60001.5
$ws.Cells.Item(73, 11).Value = 5976
$ws.Cells.Item(73, 12).Value = 180004.5
$ws.Cells.Item(73, 13).Value = -5040
$ws.Cells.Item(73, 14).Value = -181876.5

# Row 86
$ws.Cells.Item(86, 8).Value = 3117.348
$ws.Cells.Item(86, 9).Value = 1268.7273
$ws.Cells.Item(86, 10).Value = 4811.9165
$ws.Cells.Item(86, 11).Value = 1268.7273
$ws.Cells.Item(86, 12).Value = 4811.9165
$ws.Cells.Item(86, 13).Value = -145.7273
$ws.Cells.Item(86, 14).Value = -7057.9165

# Row 89
$ws.Cells.Item(89, 8).Value = 3117.348
$ws.Cells.Item(89, 9).Value = 1268.7273
$ws.Cells.Item(89, 10).Value = 4811.9165
$ws.Cells.Item(89, 11).Value = 6343.636500000001
$ws.Cells.Item(89, 12).Value = 24059.5825
$ws.Cells.Item(89, 13).Value = -727.6365000000005
$ws.Cells.Item(89, 14).Value = -35291.5825

# Row 101
$ws.Cells.Item(101, 8).Value = 1624.375
$ws.Cells.Item(101, 9).Value = 1082.5
$ws.Cells.Item(101, 11).Value = 3247.5
$ws.Cells.Item(101, 13).Value = -1625.5

# Row 132
$ws.Cells.Item(132, 8).Value = 31332.861
$ws.Cells.Item(132, 9).Value = 35633.824
$ws.Cells.Item(132, 10).Value = 4745.091
$ws.Cells.Item(132, 11).Value = 106901.472
$ws.Cells.Item(132, 12).Value = 14235.273
$ws.Cells.Item(132, 13).Value = -104371.472
$ws.Cells.Item(132, 14).Value = -19295.273

# Row 138
$ws.Cells.Item(138, 8).Value = 3399.182
$ws.Cells.Item(138, 9).Value = 3104.8215
$ws.Cells.Item(138, 11).Value = 9314.4645
$ws.Cells.Item(138, 13).Value = -4174.4645

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 1367.2709
$ws.Cells.Item(2, 9).Value = 1333.3
$ws.Cells.Item(2, 10).Value = 1537.125
$ws.Cells.Item(2, 11).Value = 1333.3
$ws.Cells.Item(2, 12).Value = 1537.125
$ws.Cells.Item(2, 13).Value = -1220.3
$ws.Cells.Item(2, 14).Value = -1763.125

# Row 32
$ws.Cells.Item(32, 8).Value = 6499841.5
$ws.Cells.Item(32, 9).Value = 6949053
$ws.Cells.Item(32, 10).Value = 31199.8
$ws.Cells.Item(32, 11).Value = 6949053
$ws.Cells.Item(32, 12).Value = 31199.8
$ws.Cells.Item(32, 13).Value = -6948766
$ws.Cells.Item(32, 14).Value = -31773.8

# Row 45
$ws.Cells.Item(45, 8).Value = 6953.0713
$ws.Cells.Item(45, 9).Value = 5112
$ws.Cells.Item(45, 11).Value = 5112
$ws.Cells.Item(45, 13).Value = -4735

# Row 116
$ws.Cells.Item(116, 8).Value = 1367.2709
$ws.Cells.Item(116, 9).Value = 1333.3
$ws.Cells.Item(116, 10).Value = 1537.125
$ws.Cells.Item(116, 11).Value = 1333.3
$ws.Cells.Item(116, 12).Value = 1537.125
$ws.Cells.Item(116, 13).Value = 960.7
$ws.Cells.Item(116, 14).Value = -6125.125

# Row 122
$ws.Cells.Item(122, 8).Value = 2819.8667
$ws.Cells.Item(122, 9).Value = 1366.4445
$ws.Cells.Item(122, 11).Value = 4099.333500000001
$ws.Cells.Item(122, 13).Value = -1649.333500000001

# Row 132
$ws.Cells.Item(132, 8).Value = 3911.2954
$ws.Cells.Item(132, 9).Value = 2853.1765
$ws.Cells.Item(132, 11).Value = 8559.529500000001
$ws.Cells.Item(132, 13).Value = -6029.529500000001

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 1367.2709
$ws.Cells.Item(3, 9).Value = 1333.3
$ws.Cells.Item(3, 10).Value = 1537.125
$ws.Cells.Item(3, 11).Value = 1333.3
$ws.Cells.Item(3, 12).Value = 1537.125
$ws.Cells.Item(3, 13).Value = -1219.3
$ws.Cells.Item(3, 14).Value = -1765.125

# Row 105
$ws.Cells.Item(105, 8).Value = 3589.6875
$ws.Cells.Item(105, 9).Value = 2885.6667
$ws.Cells.Item(105, 11).Value = 2885.6667
$ws.Cells.Item(105, 13).Value = -1138.6667

# Row 107
$ws.Cells.Item(107, 8).Value = 1744.9788
$ws.Cells.Item(107, 9).Value = 1840.0488
$ws.Cells.Item(107, 11).Value = 1840.0488
$ws.Cells.Item(107, 13).Value = 79.95119999999997

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 8671.058999999999
$ws.Cells.Item(31, 9).Value = 2976
$ws.Cells.Item(31, 10).Value = 10423.385
$ws.Cells.Item(31, 11).Value = 2976
$ws.Cells.Item(31, 12).Value = 10423.385
$ws.Cells.Item(31, 13).Value = -2681
$ws.Cells.Item(31, 14).Value = -11013.385

# Row 34
$ws.Cells.Item(34, 8).Value = 8671.058999999999
$ws.Cells.Item(34, 9).Value = 2976
$ws.Cells.Item(34, 10).Value = 10423.385
$ws.Cells.Item(34, 11).Value = 2976
$ws.Cells.Item(34, 12).Value = 10423.385
$ws.Cells.Item(34, 13).Value = -2774
$ws.Cells.Item(34, 14).Value = -10827.385

# Row 122
$ws.Cells.Item(122, 8).Value = 689.75
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 14).ClearContents()

# Row 123
$ws.Cells.Item(123, 8).Value = 299999
$ws.Cells.Item(123, 10).Value = 299999
$ws.Cells.Item(123, 12).Value = 299999
$ws.Cells.Item(123, 14).Value = -309799

# Row 125
$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(125, 14).ClearContents()

# Row 134
$ws.Cells.Item(134, 8).Value = 1484.5217
$ws.Cells.Item(134, 9).Value = 1379
$ws.Cells.Item(134, 10).Value = 2997
$ws.Cells.Item(134, 11).Value = 4137
$ws.Cells.Item(134, 12).Value = 8991
$ws.Cells.Item(134, 13).Value = -1602
$ws.Cells.Item(134, 14).Value = -14061

$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Cells.Item(7, 8).Value = 2333711
$ws.Cells.Item(7, 9).Value = 400
$ws.Cells.Item(7, 11).Value = 1200
$ws.Cells.Item(7, 13).Value = -1088

# Row 57
$ws.Cells.Item(57, 8).Value = 8903
$ws.Cells.Item(57, 10).Value = 9875
$ws.Cells.Item(57, 12).Value = 29625
$ws.Cells.Item(57, 14).Value = -30743

# Row 132
$ws.Cells.Item(132, 8).Value = 2532.0386
$ws.Cells.Item(132, 9).Value = 2047.125
$ws.Cells.Item(132, 10).Value = 2747.5557
$ws.Cells.Item(132, 11).Value = 18424.125
$ws.Cells.Item(132, 12).Value = 24728.0013
$ws.Cells.Item(132, 13).Value = -15894.125
$ws.Cells.Item(132, 14).Value = -29788.0013

$ws = $wb.Worksheets.Item("GSM")
# Row 53
$ws.Cells.Item(53, 8).Value = 30000
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 10).Value = 30000
$ws.Cells.Item(53, 11).Value = 0
$ws.Cells.Item(53, 12).Value = 30000
$ws.Cells.Item(53, 13).ClearContents()
$ws.Cells.Item(53, 14).Value = -31262

# Row 93
$ws.Cells.Item(93, 8).Value = 47999.332
$ws.Cells.Item(93, 10).Value = 47999.332
$ws.Cells.Item(93, 12).Value = 47999.332
$ws.Cells.Item(93, 14).Value = -51743.332

# Row 104
$ws.Cells.Item(104, 8).Value = 100536.8
$ws.Cells.Item(104, 10).Value = 100536.8
$ws.Cells.Item(104, 12).Value = 100536.8
$ws.Cells.Item(104, 14).Value = -107524.8

# Row 140
$ws.Cells.Item(140, 8).Value = 80996.336
$ws.Cells.Item(140, 10).Value = 80996.336
$ws.Cells.Item(140, 12).Value = 80996.336
$ws.Cells.Item(140, 14).Value = -91356.336

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Cells.Item(16, 8).Value = 20548.8
$ws.Cells.Item(16, 9).Value = 20548.8
$ws.Cells.Item(16, 11).Value = 20548.8
$ws.Cells.Item(16, 13).Value = -20378.8

# Row 22
$ws.Cells.Item(22, 8).Value = 57806.277
$ws.Cells.Item(22, 9).Value = 167596
$ws.Cells.Item(22, 10).Value = 2911.4167
$ws.Cells.Item(22, 11).Value = 167596
$ws.Cells.Item(22, 12).Value = 2911.4167
$ws.Cells.Item(22, 13).Value = -167301
$ws.Cells.Item(22, 14).Value = -3501.4167

# Row 27
$ws.Cells.Item(27, 8).Value = 57806.277
$ws.Cells.Item(27, 9).Value = 167596
$ws.Cells.Item(27, 10).Value = 2911.4167
$ws.Cells.Item(27, 11).Value = 167596
$ws.Cells.Item(27, 12).Value = 2911.4167
$ws.Cells.Item(27, 13).Value = -167489
$ws.Cells.Item(27, 14).Value = -3125.4167

# Row 61
$ws.Cells.Item(61, 8).Value = 12011.479
$ws.Cells.Item(61, 9).Value = 12917.381
$ws.Cells.Item(61, 11).Value = 12917.381
$ws.Cells.Item(61, 13).Value = -12715.381

# Row 113
$ws.Cells.Item(113, 8).Value = 12011.479
$ws.Cells.Item(113, 9).Value = 12917.381
$ws.Cells.Item(113, 11).Value = 12917.381
$ws.Cells.Item(113, 13).Value = -10747.381

# Row 132
$ws.Cells.Item(132, 8).Value = 655309.9
$ws.Cells.Item(132, 9).Value = 806257.75
$ws.Cells.Item(132, 11).Value = 2418773.25
$ws.Cells.Item(132, 13).Value = -2416243.25

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Cells.Item(2, 8).Value = 250012750
$ws.Cells.Item(2, 9).Value = 16998
$ws.Cells.Item(2, 11).Value = 16998
$ws.Cells.Item(2, 13).Value = -16886

# Row 107
$ws.Cells.Item(107, 8).Value = 2218.558
$ws.Cells.Item(107, 9).Value = 669.7143
$ws.Cells.Item(107, 10).Value = 5109.7334
$ws.Cells.Item(107, 11).Value = 2009.1429
$ws.Cells.Item(107, 12).Value = 15329.2002
$ws.Cells.Item(107, 13).Value = -89.14289999999983
$ws.Cells.Item(107, 14).Value = -19169.2002

# Row 136
$ws.Cells.Item(136, 8).Value = 6930491.5
$ws.Cells.Item(136, 9).Value = 8109480
$ws.Cells.Item(136, 11).Value = 24328440
$ws.Cells.Item(136, 13).Value = -24325890
